$wb = $excel.ActiveWorkbook

# Rename Sheet2 to "issues"
$issues = $wb.Worksheets.Item("Sheet2")
$issues.Name = "issues"

# Populate the "issues" sheet with new backlog items
# (shared-string order must match: data reset.., action new fields.., back buttons,
#  action view all fields, what should happen.., defaulters, report, unique employee id check)
$issues.Range("A9").Value = "data reset after creating new action field"
$issues.Range("A10").Value = "action new fields formatting"
$issues.Range("A11").Value = "back buttons"
$issues.Range("A12").Value = "action view all fields"
$issues.Range("A13").Value = "what should happen when rest doesn't respond"
$issues.Range("A5").Value = "defaulters"
$issues.Range("A6").Value = "report"
$issues.Range("A7").Value = "unique employee id check"

# Update the selection on Sheet1 (no longer the active tab)
$sheet1 = $wb.Worksheets.Item("Sheet1")
$sheet1.Range("A8").Select()

# Set the selection on the issues sheet and make it the active tab
$issues.Range("A8").Select()
$issues.Activate()
